$wb = $excel.ActiveWorkbook

# --- Sheet1 (About) has no text content changes; shared-string index shifts ---
# are handled automatically because the cell text stays the same.

foreach ($sheetName in @("BPaFF-BITPTaP", "BPaFF-BDTPTPF")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Add new derived rows 13-17 first so new shared strings are created
    # in the same order as the target file: lignite, hard coal, onshore wind,
    # offshore wind, crude oil, heavy or residual fuel oil, municipal solid waste.
    $ws.Range("A13").Value = "lignite"

    # Rename existing fuel-type rows
    $ws.Range("A2").Value = "hard coal"
    $ws.Range("A6").Value = "onshore wind"

    $ws.Range("A14").Value = "offshore wind"
    $ws.Range("A15").Value = "crude oil"
    $ws.Range("A16").Value = "heavy or residual fuel oil"
    $ws.Range("A17").Value = "municipal solid waste"

    # Right-align the "Boolean" header cell (adds a new cell style)
    $ws.Range("B1").HorizontalAlignment = -4152

    # Formulas for the new derived rows
    $ws.Range("B13").Formula = "=B2"
    $ws.Range("B14").Formula = "=B6"
    $ws.Range("B15").Formula = "=B11"
    $ws.Range("B16").Formula = "=B11"
    $ws.Range("B17").Formula = "=B9"
}

$wb.Save()
